$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Total" row marks: correct count (B12) and corr/total summary (E12)
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 65
$ws.Range("E12").Value = "65/140"
